$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAY-22")

# Rows 11 and 12 mirror the formatting of the existing "two-line" task rows
# (row 9 is the most recent one: A/B/C/D/E/F filled, ht 28.8, D wrapped+s10, E percent+s11).
# Copying its formatting first lets the engine reuse the existing dedup'd
# styles (s="7" for dates, s="10" for wrapped text, s="11" for percents)
# instead of minting brand-new style entries.
$ws.Range("A9:F9").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A9:F9").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)   # xlPasteFormats

# Row 13 mirrors the formatting of the existing "continuation" rows that have
# no A/B/C values (row 8 / row 3: just D/E/F filled, default row height).
$ws.Range("D8:F8").Copy()
$ws.Range("D13:F13").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# Row 11
$ws.Range("A11").Value = 7
$ws.Range("B11").Value = 44691
$ws.Range("C11").Value = "RPA GSS"
$ws.Range("D11").Value = "1. Downloading and Uploading  of the OtherSales monthly task has been completed for SC1,SSC3,SSC6,SSC8,SSC9 are completed from 01-04-2022 to 30-04-2022 as requested  Mohan san, whereas the other service centers are work in progress"
$ws.Range("E11").Value = 0.8
$ws.Range("F11").Value = "WIP"

# Row 12
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = 44692
$ws.Range("C12").Value = "RPA GSS"
$ws.Range("D12").Value = "1. Supported to Clear the communication board tasks for SSC1, SSC2,SSC3,SSC6,SSC8,SSC9,SSC10 and SSC11 and also implemented the MFA authentication, tested and it is running smoothly."
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = "Completed"

# Row 13
$ws.Range("D13").Value = "2. Supported to OtherSales task to extract the Extended warranty data and uploaded for SSC10 and SS11 and it is success  manually"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = "Completed"

$ws.Range("D21").Select()
